$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.054.92"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.657.98"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'591.20"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'189.02"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.691"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  -7.42%  "
$ws.Range("D11").Value = "'55.91"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'0.0000266"
$ws.Range("E12").Value = "  -7.83%  "
$ws.Range("D13").Value = "'10.14"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "4.240.62"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "3.652.67"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'18.72"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "67.854.46"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'12.48"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").Value = "'1.09"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").Value = "'401.27"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'4.38"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "'87.28"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'2.92"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "'10.77"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "'12.40"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "'6.07"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'3.66"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("D29").Value = "'9.23"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "'31.69"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "'7.07"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'67.88"
$ws.Range("E32").Value = "  +6.08%  "
$ws.Range("D33").Value = "'12.14"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").Value = "'43.47"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'605.08"
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.115"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "'0.386"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("E40").Value = "  -14.83%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "'2.85"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("D43").Value = "'0.0419"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "'2.48"
$ws.Range("E44").Value = "  -10.18%  "
$ws.Range("D45").Value = "'3.21"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").Value = "'0.134"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").Value = "2.745.31"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'8.81"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  -12.72%  "
